$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new (blank) column before column N ("Late"), shifting
# Late / heading(Date) / Outstanding one column to the right.
$ws.Columns("N").Insert()

# The newly inserted column picks up the width of its left neighbour
# (column M, "In Advance") rather than the default column width.
$ws.Columns("N").ColumnWidth = 9.83

# Make "Repayment schedule" the active sheet/tab (it was "Transactions"
# before) and select cell R8 on it.
$ws.Activate()
$ws.Range("R8").Select()
